$p = $ppt.ActivePresentation

# --- Slide 10 ("Prediction Test-1"): TextBox 6 (shape id 7) ---
# "The training model is able to predict..." -> "The training model can predict..."
$s10 = $p.Slides.Item(10)
$shp10 = $s10.Shapes.Item(4)
$run10 = $shp10.TextFrame.TextRange.Paragraphs(2).Runs(1)
$run10.Text = "The training model can predict the feature in the prediction image as circle with an average similarity of 73.26%"

# --- Slide 11 ("Prediction Test-2"): TextBox 3 (shape id 4) ---
# "The training model is able to predict..." -> "The training model can predict..."
$s11 = $p.Slides.Item(11)
$shp11 = $s11.Shapes.Item(4)
$run11 = $shp11.TextFrame.TextRange.Paragraphs(2).Runs(1)
$run11.Text = "The training model can predict the feature in the prediction image as rectangle with an average similarity of 71.07%"

# --- Slide 13 ("Conclusion"): Content Placeholder 2 (shape id 3) ---
# Insert a new paragraph after "...parameters for effective learning." and before "For a higher ..."
$s13 = $p.Slides.Item(13)
$shp13 = $s13.Shapes.Item(2)
$para13 = $shp13.TextFrame.TextRange.Paragraphs(2)
[void]$para13.InsertAfter("`rThe Image Dimensions influence the Potential Radius")
